$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 291; existing rows 291-310 shift down to 292-311.
$ws.Rows.Item(291).Insert()

# Populate the new row 291 with the new record (same constant columns as
# surrounding rows, new date/price data per the diff).
$ws.Cells.Item(291, 1).Value = 10
$ws.Cells.Item(291, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(291, 3).Value = "La Araucanía"
$ws.Cells.Item(291, 4).Value = 44714
$ws.Cells.Item(291, 5).Value = 9
$ws.Cells.Item(291, 6).Value = 100114013
$ws.Cells.Item(291, 7).Value = "Zanahoria"
$ws.Cells.Item(291, 8).Value = "Sin especificar"
$ws.Cells.Item(291, 9).Value = "Primera"
$ws.Cells.Item(291, 10).Value = 380
$ws.Cells.Item(291, 11).Value = 6000
$ws.Cells.Item(291, 12).Value = 6000
$ws.Cells.Item(291, 13).Value = 6000
$ws.Cells.Item(291, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(291, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(291, 16).Value = 240
$ws.Cells.Item(291, 17).Value = 25
$ws.Cells.Item(291, 18).Value = "Hortaliza"

# Keep the date cell formatted the same as the rest of column D.
$ws.Cells.Item(291, 4).NumberFormat = $ws.Cells.Item(292, 4).NumberFormat
